$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.047.04'
$ws.Range("E2").Value = '  +7.19%  '

$ws.Range("D3").Value = '3.020.15'
$ws.Range("E3").Value = '  +4.26%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.72'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.31'
$ws.Range("E6").Value = '  +8.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.21%  '

$ws.Range("D8").Value = '3.015.89'
$ws.Range("E8").Value = '  +4.14%  '

$ws.Range("E9").Value = '  +2.50%  '

$ws.Range("E10").Value = '  -0.16%  '

$ws.Range("E11").Value = '  +4.93%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.452'
$ws.Range("E12").Value = '  +4.79%  '

$ws.Range("E13").Value = '  +3.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.38'
$ws.Range("E14").Value = '  +7.96%  '

$ws.Range("E15").Value = '  +0.72%  '

$ws.Range("D16").Value = '65.937.09'
$ws.Range("E16").Value = '  +7.03%  '

$ws.Range("D17").Value = '3.520.58'
$ws.Range("E17").Value = '  +4.25%  '

$ws.Range("E18").Value = '  +6.34%  '

$ws.Range("D19").Value = '3.021.21'
$ws.Range("E19").Value = '  +4.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '460.95'
$ws.Range("E20").Value = '  +6.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.81'
$ws.Range("E21").Value = '  +5.61%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.687'
$ws.Range("E22").Value = '  +4.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.38'
$ws.Range("E23").Value = '  +8.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.13'
$ws.Range("E24").Value = '  +3.50%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.58'
$ws.Range("E25").Value = '  +5.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.26'
$ws.Range("E26").Value = '  +12.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.87'
$ws.Range("E27").Value = '  +9.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.42'
$ws.Range("E29").Value = '  +18.59%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.93'
$ws.Range("E30").Value = '  +13.31%  '

$ws.Range("E31").Value = '  -2.34%  '

$ws.Range("E32").Value = '  +4.22%  '

$ws.Range("E33").Value = '  +5.30%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.04'
$ws.Range("E34").Value = '  +5.87%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.24%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.993'
$ws.Range("E36").Value = '  +3.54%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.78'
$ws.Range("E37").Value = '  +7.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.18'
$ws.Range("E38").Value = '  +12.70%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.02'
$ws.Range("E39").Value = '  +7.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.48'
$ws.Range("E40").Value = '  +1.30%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '44.99'
$ws.Range("E41").Value = '  +13.55%  '

$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.122'
$ws.Range("E42").Value = '  +7.43%  '

$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.302'
$ws.Range("E43").Value = '  +13.47%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.47'
$ws.Range("E44").Value = '  +3.11%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '389.22'

$ws.Range("D46").Value = '2.801.37'
$ws.Range("E46").Value = '  +4.06%  '

$ws.Range("E47").Value = '  +5.26%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '135.04'
$ws.Range("E48").Value = '  +1.61%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.71'
$ws.Range("E50").Value = '  +9.79%  '

$ws.Range("E51").Value = '  +3.78%  '
